$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2231376440773829
$ws.Cells.Item(2, 3).Value = 0.04889604927413416
$ws.Cells.Item(2, 5).Value = 0.1788266674977166
$ws.Cells.Item(2, 6).Value = 0.4443680307746121
$ws.Cells.Item(2, 7).Value = 0.3111141937561257
$ws.Cells.Item(2, 8).Value = 0.4945810236140389
$ws.Cells.Item(2, 9).Value = 0.3793677714397212
$ws.Cells.Item(2, 11).Value = 0.2306778499828965
$ws.Cells.Item(2, 13).Value = 0.202850923048203
$ws.Cells.Item(2, 14).Value = 1.19242331793113
$ws.Cells.Item(2, 15).Value = 1.541002427154922
$ws.Cells.Item(3, 2).Value = 0.1948131132054129
$ws.Cells.Item(3, 3).Value = 0.04578998727642869
$ws.Cells.Item(3, 5).Value = 0.1673073147037982
$ws.Cells.Item(3, 6).Value = 0.3878228170618172
$ws.Cells.Item(3, 7).Value = 0.3135636137168873
$ws.Cells.Item(3, 8).Value = 0.4986741276632785
$ws.Cells.Item(3, 9).Value = 0.3838434522605212
$ws.Cells.Item(3, 11).Value = 0.2016878197841692
$ws.Cells.Item(3, 13).Value = 0.1813066637136913
$ws.Cells.Item(3, 14).Value = 1.203389232103483
$ws.Cells.Item(3, 15).Value = 1.554568244127665
$ws.Cells.Item(4, 2).Value = 0.1773689323459848
$ws.Cells.Item(4, 3).Value = 0.04386811979645699
$ws.Cells.Item(4, 5).Value = 0.1603589720155654
$ws.Cells.Item(4, 6).Value = 0.3531389305169483
$ws.Cells.Item(4, 7).Value = 0.315298622954252
$ws.Cells.Item(4, 8).Value = 0.5013908483636484
$ws.Cells.Item(4, 9).Value = 0.3867976588554782
$ws.Cells.Item(4, 11).Value = 0.18380753118538
$ws.Cells.Item(4, 13).Value = 0.168125745538859
$ws.Cells.Item(4, 14).Value = 1.21053067326293
$ws.Cells.Item(4, 15).Value = 1.56380694433102
$ws.Cells.Item(5, 2).Value = 0.1702475059569792
$ws.Cells.Item(5, 3).Value = 0.04308128116201715
$ws.Cells.Item(5, 5).Value = 0.1575586467702692
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.3160636909006627
$ws.Cells.Item(5, 8).Value = 0.5025491472004049
$ws.Cells.Item(5, 9).Value = 0.3880533381300157
$ws.Cells.Item(5, 11).Value = 0.1765014538777336
$ws.Cells.Item(5, 13).Value = 0.1627663993555899
$ws.Cells.Item(5, 14).Value = 1.21354360313213
$ws.Cells.Item(5, 15).Value = 1.567800351858097
$ws.Cells.Item(6, 2).Value = 0.1690642412153807
$ws.Cells.Item(6, 3).Value = 0.04295040728445798
$ws.Cells.Item(6, 5).Value = 0.1570955338275866
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.3161942333241186
$ws.Cells.Item(6, 8).Value = 0.5027445759211915
$ws.Cells.Item(6, 9).Value = 0.3882649718562483
$ws.Cells.Item(6, 11).Value = 0.1752871086026033
$ws.Cells.Item(6, 13).Value = 0.1618772102325039
$ws.Cells.Item(6, 14).Value = 1.214050104496266
$ws.Cells.Item(6, 15).Value = 1.568477256626522
$ws.Cells.Item(7, 2).Value = 0.1772729413905836
$ws.Cells.Item(7, 3).Value = 0.04385752298590972
$ws.Cells.Item(7, 5).Value = 0.1603210797952173
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.315308706019735
$ws.Cells.Item(7, 8).Value = 0.5014062621843181
$ws.Cells.Item(7, 9).Value = 0.3868143836030278
$ws.Cells.Item(7, 11).Value = 0.1837090780645951
$ws.Cells.Item(7, 13).Value = 0.1680534188892508
$ws.Cells.Item(7, 14).Value = 1.210570890642249
$ws.Cells.Item(7, 15).Value = 1.563859875477362
$ws.Cells.Item(8, 2).Value = 0.213382595068083
$ws.Cells.Item(8, 3).Value = 0.04782816160599168
$ws.Cells.Item(8, 5).Value = 0.1748288443796184
$ws.Cells.Item(8, 6).Value = 0.4248636149813478
$ws.Cells.Item(8, 7).Value = 0.3119107599396358
$ws.Cells.Item(8, 8).Value = 0.4959501006555698
$ws.Cells.Item(8, 9).Value = 0.3808681974101571
$ws.Cells.Item(8, 11).Value = 0.2206990292187356
$ws.Cells.Item(8, 13).Value = 0.1954126776306566
$ws.Cells.Item(8, 14).Value = 1.196119652352536
$ws.Cells.Item(8, 15).Value = 1.545491165605753
$ws.Cells.Item(9, 2).Value = 0.2837559710762037
$ws.Cells.Item(9, 3).Value = 0.05549614104037914
$ws.Cells.Item(9, 5).Value = 0.2042753025011095
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.3070837628436891
$ws.Cells.Item(9, 8).Value = 0.4868641692390341
$ws.Cells.Item(9, 9).Value = 0.370843814341514
$ws.Cells.Item(9, 11).Value = 0.292581514181677
$ws.Cells.Item(9, 13).Value = 0.2494389774421819
$ws.Cells.Item(9, 14).Value = 1.17101789557578
$ws.Cells.Item(9, 15).Value = 1.516688270905007
$ws.Cells.Item(10, 2).Value = 0.3351724137824874
$ws.Cells.Item(10, 3).Value = 0.06105605197505781
$ws.Cells.Item(10, 5).Value = 0.2265309612139887
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.3046613552039048
$ws.Cells.Item(10, 8).Value = 0.4811705301362394
$ws.Cells.Item(10, 9).Value = 0.3644770313019627
$ws.Cells.Item(10, 11).Value = 0.3449756374940307
$ws.Cells.Item(10, 13).Value = 0.2893636312822849
$ws.Cells.Item(10, 14).Value = 1.154544889502063
$ws.Cells.Item(10, 15).Value = 1.499932977249387
$ws.Cells.Item(11, 2).Value = 0.3584968044839343
$ws.Cells.Item(11, 3).Value = 0.06356908509640391
$ws.Cells.Item(11, 5).Value = 0.2367936162133049
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.3038043828044508
$ws.Cells.Item(11, 8).Value = 0.4787931427263246
$ws.Cells.Item(11, 9).Value = 0.3617975124120214
$ws.Cells.Item(11, 11).Value = 0.3687166234080053
$ws.Cells.Item(11, 13).Value = 0.3075775304692669
$ws.Cells.Item(11, 14).Value = 1.147477642717909
$ws.Cells.Item(11, 15).Value = 1.493268644957723
$ws.Cells.Item(12, 2).Value = 0.367319321488452
$ws.Cells.Item(12, 3).Value = 0.06451833915025418
$ws.Cells.Item(12, 5).Value = 0.2406999177760554
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.3035151744121194
$ws.Cells.Item(12, 8).Value = 0.4779234427408383
$ws.Cells.Item(12, 9).Value = 0.3608140370905488
$ws.Cells.Item(12, 11).Value = 0.3776928675058002
$ws.Cells.Item(12, 13).Value = 0.3144821079294715
$ws.Cells.Item(12, 14).Value = 1.144862729199254
$ws.Cells.Item(12, 15).Value = 1.490882861965417
$ws.Cells.Item(13, 2).Value = 0.365419684653375
$ws.Cells.Item(13, 3).Value = 0.06431400678539489
$ws.Cells.Item(13, 5).Value = 0.2398577314522043
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.3035758890133238
$ws.Cells.Item(13, 8).Value = 0.4781093890864057
$ws.Cells.Item(13, 9).Value = 0.3610244579760362
$ws.Cells.Item(13, 11).Value = 0.3757603016329369
$ws.Cells.Item(13, 13).Value = 0.3129947575717651
$ws.Cells.Item(13, 14).Value = 1.145423171554199
$ws.Cells.Item(13, 15).Value = 1.491390549759274
$ws.Cells.Item(14, 2).Value = 0.3592228402818307
$ws.Cells.Item(14, 3).Value = 0.06364722871110473
$ws.Cells.Item(14, 5).Value = 0.2371145871768903
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.3037798813730106
$ws.Cells.Item(14, 8).Value = 0.4787209795169929
$ws.Cells.Item(14, 9).Value = 0.3617159758321229
$ws.Cells.Item(14, 11).Value = 0.3694553870460595
$ws.Cells.Item(14, 13).Value = 0.3081454271284585
$ws.Cells.Item(14, 14).Value = 1.147261283767669
$ws.Cells.Item(14, 15).Value = 1.493069601571648
$ws.Cells.Item(15, 2).Value = 0.3554257856096115
$ws.Cells.Item(15, 3).Value = 0.06323849690981831
$ws.Cells.Item(15, 5).Value = 0.2354369484894363
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.3039094331138088
$ws.Cells.Item(15, 8).Value = 0.4790995762616959
$ws.Cells.Item(15, 9).Value = 0.3621436145825498
$ws.Cells.Item(15, 11).Value = 0.365591612226325
$ws.Cells.Item(15, 13).Value = 0.3051760287589929
$ws.Cells.Item(15, 14).Value = 1.148395163096986
$ws.Cells.Item(15, 15).Value = 1.494116026039023
$ws.Cells.Item(16, 2).Value = 0.3336467511101375
$ws.Cells.Item(16, 3).Value = 0.06089148931559407
$ws.Cells.Item(16, 5).Value = 0.2258630711663585
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.3047222969981576
$ws.Cells.Item(16, 8).Value = 0.4813301760893722
$ws.Cells.Item(16, 9).Value = 0.3646565085588058
$ws.Cells.Item(16, 11).Value = 0.3434221845077161
$ws.Cells.Item(16, 13).Value = 0.2881743470690239
$ws.Cells.Item(16, 14).Value = 1.155015329560953
$ws.Cells.Item(16, 15).Value = 1.500387790578657
$ws.Cells.Item(17, 2).Value = 0.3202689286002567
$ws.Cells.Item(17, 3).Value = 0.05944749340980593
$ws.Cells.Item(17, 5).Value = 0.2200253726024215
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.3052837674781799
$ws.Cells.Item(17, 8).Value = 0.4827530359709442
$ws.Cells.Item(17, 9).Value = 0.366253627194979
$ws.Cells.Item(17, 11).Value = 0.3297976713468245
$ws.Cells.Item(17, 13).Value = 0.2777576037814242
$ws.Cells.Item(17, 14).Value = 1.159185787179776
$ws.Cells.Item(17, 15).Value = 1.504480709409052
$ws.Cells.Item(18, 2).Value = 0.312568254002656
$ws.Cells.Item(18, 3).Value = 0.05861542335243541
$ws.Cells.Item(18, 5).Value = 0.2166807049992627
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.305629765677395
$ws.Cells.Item(18, 8).Value = 0.4835914466087416
$ws.Cells.Item(18, 9).Value = 0.3671926485228951
$ws.Cells.Item(18, 11).Value = 0.3219524553597353
$ws.Cells.Item(18, 13).Value = 0.2717710669312936
$ws.Cells.Item(18, 14).Value = 1.161624660367085
$ws.Cells.Item(18, 15).Value = 1.506924986650148
$ws.Cells.Item(19, 2).Value = 0.3099599060423941
$ws.Cells.Item(19, 3).Value = 0.05833343881211306
$ws.Cells.Item(19, 5).Value = 0.215550490043988
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.3057508720704618
$ws.Cells.Item(19, 8).Value = 0.4838787569171288
$ws.Cells.Item(19, 9).Value = 0.3675140879703989
$ws.Cells.Item(19, 11).Value = 0.3192947121567329
$ws.Cells.Item(19, 13).Value = 0.2697449728828758
$ws.Cells.Item(19, 14).Value = 1.162457314928115
$ws.Cells.Item(19, 15).Value = 1.507768054221103
$ws.Cells.Item(20, 2).Value = 0.3216936564420791
$ws.Cells.Item(20, 3).Value = 0.0596013670874953
$ws.Cells.Item(20, 5).Value = 0.2206454567504395
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.3052216112848711
$ws.Cells.Item(20, 8).Value = 0.4825994983173132
$ws.Cells.Item(20, 9).Value = 0.3660814995514059
$ws.Cells.Item(20, 11).Value = 0.3312489340974309
$ws.Cells.Item(20, 13).Value = 0.2788659776781728
$ws.Cells.Item(20, 14).Value = 1.158737681024899
$ws.Cells.Item(20, 15).Value = 1.50403568123393
$ws.Cells.Item(21, 2).Value = 0.3610432773666901
$ws.Cells.Item(21, 3).Value = 0.06384314248690259
$ws.Cells.Item(21, 5).Value = 0.2379197696286894
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.3037190049861564
$ws.Cells.Item(21, 8).Value = 0.4785405111969467
$ws.Cells.Item(21, 9).Value = 0.3615120132675571
$ws.Cells.Item(21, 11).Value = 0.3713076761418961
$ws.Cells.Item(21, 13).Value = 0.3095695938740519
$ws.Cells.Item(21, 14).Value = 1.14671972213106
$ws.Cells.Item(21, 15).Value = 1.4925726809455
$ws.Cells.Item(22, 2).Value = 0.3867024337420446
$ws.Cells.Item(22, 3).Value = 0.06660151054551022
$ws.Cells.Item(22, 5).Value = 0.2493265065342172
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.3029428057241432
$ws.Cells.Item(22, 8).Value = 0.4760658767717416
$ws.Cells.Item(22, 9).Value = 0.3587074631385967
$ws.Cells.Item(22, 11).Value = 0.3974068126132977
$ws.Cells.Item(22, 13).Value = 0.3296790853531135
$ws.Cells.Item(22, 14).Value = 1.13922258493956
$ws.Cells.Item(22, 15).Value = 1.485884533599943
$ws.Cells.Item(23, 2).Value = 0.373013155450991
$ws.Cells.Item(23, 3).Value = 0.06513060314337338
$ws.Cells.Item(23, 5).Value = 0.243227764318533
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.3033382175659867
$ws.Cells.Item(23, 8).Value = 0.477370340966047
$ws.Cells.Item(23, 9).Value = 0.3601876529735364
$ws.Cells.Item(23, 11).Value = 0.3834848495687027
$ws.Cells.Item(23, 13).Value = 0.3189423805678757
$ws.Cells.Item(23, 14).Value = 1.143191260899862
$ws.Cells.Item(23, 15).Value = 1.48938055152928
$ws.Cells.Item(24, 2).Value = 0.3210495668033957
$ws.Cells.Item(24, 3).Value = 0.0595318067047117
$ws.Cells.Item(24, 5).Value = 0.2203650808892874
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.3052496398357363
$ws.Cells.Item(24, 8).Value = 0.4826688491771876
$ws.Cells.Item(24, 9).Value = 0.3661592536195464
$ws.Cells.Item(24, 11).Value = 0.3305928564869873
$ws.Cells.Item(24, 13).Value = 0.2783648750192427
$ws.Cells.Item(24, 14).Value = 1.15894014141335
$ws.Cells.Item(24, 15).Value = 1.504236594389454
$ws.Cells.Item(25, 2).Value = 0.2647670335232135
$ws.Cells.Item(25, 3).Value = 0.05343458928889788
$ws.Cells.Item(25, 5).Value = 0.1962011253456737
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.3081925321140275
$ws.Cells.Item(25, 8).Value = 0.4891496054573778
$ws.Cells.Item(25, 9).Value = 0.3733804466501489
$ws.Cells.Item(25, 11).Value = 0.2732074337890822
$ws.Cells.Item(25, 13).Value = 0.2347829843037488
$ws.Cells.Item(25, 14).Value = 1.177462472837103
$ws.Cells.Item(25, 15).Value = 1.523706864350629
